$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28
$ws.Range("D28").Value = "[임피던스 제어(1)] 임피던스 제어란?"
$ws.Range("E28").Value = "https://ropiens.tistory.com/94"

# Row 32
$ws.Range("D32").Value = "회귀분석의 가정과 한계 극복 방법 (nc 단비 블로그 퍼옴)"
$ws.Range("E32").Value = "https://dodonam.tistory.com/251"

# Row 37
$ws.Range("D37").Value = "[Paper Review] A Simple Framework for Contrastive Learning of Visual Representations"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1448&mod=document&pageid=1"

# Row 39
$ws.Range("D39").Value = "Probability concepts explained: Marginalisation"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/Probability-concepts-explained-Marginalisation-1"

# Row 50
$ws.Range("D50").Value = "Miura-ori [origami]"
$ws.Range("E50").Value = "http://incredible.egloos.com/7512931"

# Row 51
$ws.Range("D51").Value = "[세이버메트릭스] 인플레이 타구가 안타가 될 확률, BABIP"
$ws.Range("E51").Value = "https://bskyvision.com/1135"
